$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "WolframAlpha Value" column (D1:D22) contents.
# D1 is the header cell, D2:D22 hold the numeric WolframAlpha values.
$ws.Range("D1:D22").Value = $null

# Update the selection to reflect the now-active column D range.
$ws.Range("D1:D22").Select()
